$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$values = @{
    2  = "2025-10-29T13:14:08.772524+00:00"
    3  = "2025-10-29T13:14:11.387724+00:00"
    4  = "2025-10-29T13:14:11.387741+00:00"
    5  = "2025-10-29T13:14:11.387749+00:00"
    6  = "2025-10-29T13:14:11.387757+00:00"
    7  = "2025-10-29T13:14:11.387765+00:00"
    8  = "2025-10-29T13:14:13.474366+00:00"
    9  = "2025-10-29T13:14:13.474387+00:00"
    10 = "2025-10-29T13:14:15.942585+00:00"
    11 = "2025-10-29T13:14:18.408390+00:00"
    12 = "2025-10-29T13:14:18.408406+00:00"
    13 = "2025-10-29T13:14:20.502856+00:00"
    14 = "2025-10-29T13:14:20.502872+00:00"
    15 = "2025-10-29T13:14:20.502880+00:00"
    16 = "2025-10-29T13:14:20.502887+00:00"
    17 = "2025-10-29T13:14:27.107801+00:00"
    18 = "2025-10-29T13:14:29.151993+00:00"
    19 = "2025-10-29T13:14:31.281310+00:00"
    20 = "2025-10-29T13:14:33.679217+00:00"
    21 = "2025-10-29T13:14:33.679236+00:00"
    22 = "2025-10-29T13:14:33.679244+00:00"
    23 = "2025-10-29T13:14:35.717396+00:00"
    24 = "2025-10-29T13:14:35.717412+00:00"
    25 = "2025-10-29T13:14:35.717420+00:00"
    26 = "2025-10-29T13:14:35.717427+00:00"
    27 = "2025-10-29T13:14:35.717434+00:00"
    28 = "2025-10-29T13:14:42.714248+00:00"
    29 = "2025-10-29T13:14:42.714282+00:00"
    30 = "2025-10-29T13:14:42.714302+00:00"
    31 = "2025-10-29T13:14:42.714319+00:00"
    32 = "2025-10-29T13:14:45.016888+00:00"
    33 = "2025-10-29T13:14:45.016918+00:00"
    34 = "2025-10-29T13:14:45.016936+00:00"
    35 = "2025-10-29T13:14:47.420790+00:00"
    36 = "2025-10-29T13:14:47.420816+00:00"
    37 = "2025-10-29T13:14:47.420833+00:00"
    38 = "2025-10-29T13:14:47.420848+00:00"
    39 = "2025-10-29T13:14:47.420865+00:00"
    40 = "2025-10-29T13:14:47.420880+00:00"
    41 = "2025-10-29T13:14:47.420897+00:00"
    42 = "2025-10-29T13:14:47.420912+00:00"
    43 = "2025-10-29T13:14:47.420926+00:00"
    44 = "2025-10-29T13:14:49.719519+00:00"
    45 = "2025-10-29T13:14:49.719547+00:00"
    46 = "2025-10-29T13:14:54.799374+00:00"
    47 = "2025-10-29T13:14:57.509386+00:00"
    48 = "2025-10-29T13:14:57.509414+00:00"
    49 = "2025-10-29T13:14:57.509431+00:00"
    50 = "2025-10-29T13:14:57.509447+00:00"
}

foreach ($row in $values.Keys) {
    $ws.Range("K$row").Value = $values[$row]
}
